# Fruta / hortaliza, semanal
# Insert two new weekly-report rows (new row 15 and 16) above the existing
# data, pushing the previous rows 15-25 down to rows 17-27, then populate
# the two freshly inserted rows with their own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Shift existing rows 15:25 down by two rows (creates two blank rows at 15:16,
# inheriting the formatting - including the date number format on column D -
# from the row immediately below, same as Excel's native Insert behaviour).
$ws.Rows("15:16").Insert()

# --- New row 15 ---
$ws.Cells.Item(15, 1).Value = 6
$ws.Cells.Item(15, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(15, 3).Value = "Metropolitana"
$ws.Cells.Item(15, 4).Value = 44586
$ws.Cells.Item(15, 5).Value = 13
$ws.Cells.Item(15, 6).Value = "Fruta"
$ws.Cells.Item(15, 7).Value = 100102
$ws.Cells.Item(15, 8).Value = "Cítricos"
$ws.Cells.Item(15, 9).Value = 100102006
$ws.Cells.Item(15, 10).Value = "Pomelo"
$ws.Cells.Item(15, 11).Value = "Start Ruby"
$ws.Cells.Item(15, 12).Value = "Primera"
$ws.Cells.Item(15, 13).Value = 20
$ws.Cells.Item(15, 14).Value = 180000
$ws.Cells.Item(15, 15).Value = 180000
$ws.Cells.Item(15, 16).Value = 180000
$ws.Cells.Item(15, 17).Value = "$/bins (350 kilos)"
$ws.Cells.Item(15, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(15, 19).Value = 514
$ws.Cells.Item(15, 20).Value = 350

# --- New row 16 ---
$ws.Cells.Item(16, 1).Value = 6
$ws.Cells.Item(16, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(16, 3).Value = "Metropolitana"
$ws.Cells.Item(16, 4).Value = 44586
$ws.Cells.Item(16, 5).Value = 13
$ws.Cells.Item(16, 6).Value = "Fruta"
$ws.Cells.Item(16, 7).Value = 100102
$ws.Cells.Item(16, 8).Value = "Cítricos"
$ws.Cells.Item(16, 9).Value = 100102006
$ws.Cells.Item(16, 10).Value = "Pomelo"
$ws.Cells.Item(16, 11).Value = "Start Ruby"
$ws.Cells.Item(16, 12).Value = "Segunda"
$ws.Cells.Item(16, 13).Value = 12
$ws.Cells.Item(16, 14).Value = 140000
$ws.Cells.Item(16, 15).Value = 140000
$ws.Cells.Item(16, 16).Value = 140000
$ws.Cells.Item(16, 17).Value = "$/bins (350 kilos)"
$ws.Cells.Item(16, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(16, 19).Value = 400
$ws.Cells.Item(16, 20).Value = 350
